$d = $word.ActiveDocument

# Update the date/title line
$d.Content.Find.Execute("2024-03-10 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-11 Monday", 2)

# Update each cell in the practice table (20 rows x 5 columns, row-major order)
$values = @(
    "28+69=",
    "5+47=",
    "73-58=",
    "28+56=",
    "51-17=",
    "70-24=",
    "76-37=",
    "66+28=",
    "55-26=",
    "31-28=",
    "77+15=",
    "40-29=",
    "60-52=",
    "66-17=",
    "6+78=",
    "78+4=",
    "75-18=",
    "7+77=",
    "73-44=",
    "15+16=",
    "43+48=",
    "9+22=",
    "4+28=",
    "50-32=",
    "23+9=",
    "13+68=",
    "93-88=",
    "5+39=",
    "61-35=",
    "97-49=",
    "27+38=",
    "45+47=",
    "82-53=",
    "14+19=",
    "98-79=",
    "81-56=",
    "71-54=",
    "40-9=",
    "15+77=",
    "93-8=",
    "43+38=",
    "26+28=",
    "9+79=",
    "7+14=",
    "13+9=",
    "20-14=",
    "32+59=",
    "6+58=",
    "47+7=",
    "68+15=",
    "66+19=",
    "82-13=",
    "9+73=",
    "46+17=",
    "92-58=",
    "91-15=",
    "88+3=",
    "53+18=",
    "5+47=",
    "42+29=",
    "15+26=",
    "47+26=",
    "94-46=",
    "9+37=",
    "73-25=",
    "24-9=",
    "5+86=",
    "98-89=",
    "66-27=",
    "90-21=",
    "22+69=",
    "41-4=",
    "55-39=",
    "63-56=",
    "79+19=",
    "88+5=",
    "41-12=",
    "11-7=",
    "83-49=",
    "25+36=",
    "92-44=",
    "5+46=",
    "34+17=",
    "47+46=",
    "29+38=",
    "86-17=",
    "27+66=",
    "4+79=",
    "34+37=",
    "75+8=",
    "8+44=",
    "45+48=",
    "47+27=",
    "52-8=",
    "37+35=",
    "96-77=",
    "7+8=",
    "98-79=",
    "95-29=",
    "77-59="
)

$t = $d.Tables.Item(1)
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = [int][math]::Floor($i / 5) + 1
    $col = ($i % 5) + 1
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $values[$i]
}